$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B6").Value = "Sub23"
$ws.Range("C6").Value = "caballeros"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Minigutti, Ignacio"
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = "'"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 100
